# Matriz de adyacencia - fill in adjacency values for actors 18-21 (rows 20-23)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Matriz de adyacencia")

$row20 = @(0,1,0,0,0,1,1,0,0,1,0,0,0,0,0,1,1,0,0,1,0,0,0,1,0,1,0,1,1,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0)
for ($i = 0; $i -lt $row20.Length; $i++) {
    $ws.Cells.Item(20, 3 + $i).Value = $row20[$i]
}

$row21 = @(0,1,0,0,0,0,0,0,1,1,0,0,0,0,0,1,0,1,0,0,1,0,1,0,1,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1)
for ($i = 0; $i -lt $row21.Length; $i++) {
    $ws.Cells.Item(21, 3 + $i).Value = $row21[$i]
}

$row22 = @(0,0,0,0,0,1,1,0,0,1,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $row22.Length; $i++) {
    $ws.Cells.Item(22, 3 + $i).Value = $row22[$i]
}

$row23 = @(1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0)
for ($i = 0; $i -lt $row23.Length; $i++) {
    $ws.Cells.Item(23, 3 + $i).Value = $row23[$i]
}

# W22 loses its special "no right border" column-separator style once re-entered;
# reuse the plain bordered style already used by neighbouring data cells (e.g. C22)
# via a format-only paste so no new style entry gets created.
$ws.Range("C22").Copy() | Out-Null
$ws.Range("W22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Leave the selection where the user ended up after entering the data
$ws.Activate()
$ws.Range("AY30").Select() | Out-Null
